$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '43.121.01'
$ws.Range('E2').Value = '  +0.25%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '2.321.42'
$ws.Range('E3').Value = '  +0.96%  '

$ws.Range('E4').Value = '  +0.05%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '303.34'
$ws.Range('E5').Value = '  +0.56%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '99.74'
$ws.Range('E6').Value = '  +0.77%  '

$ws.Range('E7').Value = '  +0.29%  '

$ws.Range('E8').Value = '  +0.05%  '

$ws.Range('E9').Value = '  +2.56%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '36.24'
$ws.Range('E10').Value = '  +6.13%  '

$ws.Range('E11').Value = '  -0.42%  '

$ws.Range('E12').Value = '  -0.97%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '17.78'
$ws.Range('E13').Value = '  -0.93%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '6.94'
$ws.Range('E14').Value = '  +2.07%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '2.683.56'
$ws.Range('E15').Value = '  +1.02%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '2.308.55'
$ws.Range('E16').Value = '  +1.93%  '

$ws.Range('E17').Value = '  -1.29%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '43.045.15'
$ws.Range('E18').Value = '  +0.28%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '12.99'
$ws.Range('E19').Value = '  +4.57%  '

$ws.Range('E20').Value = '  +2.08%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '0.0₃0912'
$ws.Range('E21').Value = '  +0.65%  '

$ws.Range('E22').Value = '  +0.80%  '

$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '240.92'
$ws.Range('E23').Value = '  +1.75%  '

$ws.Range('E24').Value = '  -1.99%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '2.46'
$ws.Range('E25').Value = '  +0.07%  '

$ws.Range('E26').Value = '  -0.09%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '25.51'
$ws.Range('E27').Value = '  +3.11%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '168.80'
$ws.Range('E28').Value = '  +0.73%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '34.25'
$ws.Range('E29').Value = '  +1.11%  '

$ws.Range('E30').Value = '  +0.64%  '

$ws.Range('E31').Value = '  -1.79%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.96'
$ws.Range('E32').Value = '  +9.11%  '

$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '5.16'
$ws.Range('E33').Value = '  +2.33%  '

$ws.Range('B34').Value = 'FirstDigitalUSD'
$ws.Range('C34').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('E34').Value = '  -0.03%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '17.81'
$ws.Range('E35').Value = '  +5.19%  '

$ws.Range('E36').Value = '  -1.53%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.0700'
$ws.Range('E37').Value = '  +1.29%  '

$ws.Range('E38').Value = '  +1.86%  '

$ws.Range('E39').Value = '  +0.01%  '

$ws.Range('E40').Value = '  -0.56%  '

$ws.Range('E41').Value = '  +0.18%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '1.994.60'
$ws.Range('E42').Value = '  +0.14%  '

$ws.Range('E43').Value = '  +1.71%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '10.17'
$ws.Range('E45').Value = '  +1.35%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '17.53'
$ws.Range('E46').Value = '  -0.26%  '

$ws.Range('E47').Value = '  +0.22%  '

$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '76.28'
$ws.Range('E48').Value = '  +8.88%  '

$ws.Range('E49').Value = '  -2.32%  '

$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '2.549.59'
$ws.Range('E50').Value = '  +0.97%  '

$ws.Range('E51').Value = '  +12.33%  '
